$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.823.86"
$ws.Range("E2").Value = "  +1.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.880.01"
$ws.Range("E3").Value = "  +1.88%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.57"
$ws.Range("E5").Value = "  +3.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4738"
$ws.Range("E7").Value = "  +6.68%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3979"
$ws.Range("E8").Value = "  +4.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.56"
$ws.Range("E9").Value = "  +0.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08044"
$ws.Range("E10").Value = "  +2.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.026"
$ws.Range("E11").Value = "  +1.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.01"
$ws.Range("E12").Value = "  +3.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.896.07"
$ws.Range("E13").Value = "  +3.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.970"
$ws.Range("E14").Value = "  +2.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.190"
$ws.Range("E15").Value = "  +1.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001053"
$ws.Range("E17").Value = "  +2.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "87.23"
$ws.Range("E18").Value = "  +2.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06632"
$ws.Range("E19").Value = "  +2.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.35"
$ws.Range("E20").Value = "  +2.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.976.29"
$ws.Range("E22").Value = "  +2.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.510"
$ws.Range("E23").Value = "  +1.33%  "
$ws.Range("E24").Value = "  +3.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.302"
$ws.Range("E25").Value = "  +2.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.138.61"
$ws.Range("E26").Value = "  +4.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.24"
$ws.Range("E27").Value = "  +3.82%  "
$ws.Range("E28").Value = "  +5.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.108"
$ws.Range("E29").Value = "  +3.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.609"
$ws.Range("E30").Value = "  +2.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.85"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9789"
$ws.Range("E32").Value = "  +6.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09574"
$ws.Range("E33").Value = "  +3.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.468"
$ws.Range("E34").Value = "  +1.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.638"
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.322"
$ws.Range("E36").Value = "  +1.94%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06119"
$ws.Range("E37").Value = "  +3.14%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02264"
$ws.Range("E38").Value = "  +2.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.235"
$ws.Range("E39").Value = "  +3.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.236"
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6037"
$ws.Range("E41").Value = "  +2.88%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("E43").Value = "  +4.03%  "
$ws.Range("E44").Value = "  +1.27%  "
$ws.Range("E45").Value = "  +2.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.246"
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.27"
$ws.Range("E47").Value = "  +1.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.418"
$ws.Range("E48").Value = "  +1.77%  "
$ws.Range("E49").Value = "  +1.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06830"
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "113.88"
$ws.Range("E51").Value = "  +5.54%  "
